# edit.ps1 - applies the commit's changes to the deck:
#   1. Inserts a new slide at position 4 ("Protein Structure and Domains" /
#      "Receptor Binding Domain (RBD) is at amino acids 327-529") using the
#      "Title and Content" layout, pushing the former slide 4
#      ("Change in V2 along RBD-Down Proteins") and everything after it down
#      by one.
#   2. Reflows / raises the title box on slide 1 (the title slide) so it
#      sits higher and taller than before.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Insert the new "Protein Structure and Domains" slide at index 4.
# ---------------------------------------------------------------------
$layout = $p.SlideMaster.CustomLayouts.Item(2)
$newSlide = $p.Slides.AddSlide(4, $layout)

foreach ($ph in $newSlide.Shapes.Placeholders) {
    if ($ph.PlaceholderFormat.Type -eq 1) {
        $titleShape = $ph
    } elseif ($ph.PlaceholderFormat.Type -eq 7) {
        $bodyShape = $ph
    }
}

$titleShape.TextFrame.TextRange.Text = "Protein Structure and Domains"
$titleShape.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$bodyTr = $bodyShape.TextFrame.TextRange
$bodyTr.Text = "Receptor Binding Domain (RBD) is at amino acids 327-529"

$rbdStart = "Receptor Binding Domain (".Length + 1
$rbdLen = "RBD".Length
$bodyTr.Characters($rbdStart, $rbdLen).Font.Underline = -1

# ---------------------------------------------------------------------
# 2) Move/resize the title placeholder on slide 1.
#    Target EMU: off x=6600824 y=9526  ext cx=4749800 cy=3199342
#    (Left/Width are unchanged; only Top/Height move.)
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$titleOnSlide1 = $s1.Shapes.Item("Title 1")

# Point values chosen so that the engine's point->EMU conversion lands
# exactly on the target EMU (914400 EMU/in, 12700 EMU/pt).
$titleOnSlide1.Top = 0.7500797401574804
$titleOnSlide1.Height = 251.9166949133858
